$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 856
$ws1.Range("F6").Value = 660
$ws1.Range("F7").Value = 1233
$ws1.Range("F9").Value = 813
$ws1.Range("F10").Value = 693
$ws1.Range("F11").Value = 264
$ws1.Range("F13").Value = 357
$ws1.Range("F14").Value = 723
$ws1.Range("F15").Value = 944
$ws1.Range("F16").Value = 9921
$ws1.Range("F17").Value = 620
$ws1.Range("F18").Value = 48
$ws1.Range("F23").Value = 1761
$ws1.Range("F24").Value = 28
$ws1.Range("F26").Value = 488
$ws1.Range("F27").Value = 182
$ws1.Range("F28").Value = 106
$ws1.Range("F29").Value = 274
$ws1.Range("F32").Value = 68
$ws1.Range("F36").Value = 195
$ws1.Range("F37").Value = 172
$ws1.Range("F38").Value = 40
$ws1.Range("F39").Value = 92

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value = 132
$ws2.Range("F10").Value = 239
$ws2.Range("F15").Value = 61
$ws2.Range("F16").Value = 280
$ws2.Range("F21").Value = 3

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 820

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 820
$ws4.Range("F9").Value = 856
$ws4.Range("F10").Value = 660
$ws4.Range("F11").Value = 1233
$ws4.Range("F13").Value = 132
$ws4.Range("F14").Value = 813
$ws4.Range("F15").Value = 693
$ws4.Range("F16").Value = 264
$ws4.Range("F17").Value = 357
$ws4.Range("F19").Value = 944
$ws4.Range("F20").Value = 9921
$ws4.Range("F21").Value = 239
$ws4.Range("F22").Value = 620
$ws4.Range("F23").Value = 48
$ws4.Range("F26").Value = 1761
$ws4.Range("F27").Value = 28
$ws4.Range("F28").Value = 488
$ws4.Range("F29").Value = 182
$ws4.Range("F33").Value = 106
$ws4.Range("F35").Value = 61
$ws4.Range("F36").Value = 274
$ws4.Range("F39").Value = 68
$ws4.Range("F46").Value = 195
$ws4.Range("F47").Value = 172
